# "Fixed linking + IDs"
#
# The "ID" field of the MIF form is renamed to "REGID":
#  - survey!F61 (the field's `name` column) changes from "ID" to "REGID"
#    (its English/Portuguese prompt text in G61/H61 stays "ID").
#  - model sheet: the old model row for "ID" (row 20) is removed and a
#    new model row for "REGID" (integer, not a session variable) is
#    appended at the end of the block (new row 42), shifting the rows
#    that used to sit between them up by one.
#
# Also disables iterative calculation (workbook calc options).

$wb = $excel.ActiveWorkbook

# --- Disable iterative calculation -----------------------------------
$excel.Iteration = $false

# --- survey sheet: rename the "ID" field to "REGID" -------------------
$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Range("F61").Value = "REGID"

# --- model sheet: drop old "ID" row, append new "REGID" row -----------
$wsModel = $wb.Worksheets.Item("model")

# Remove the old row 20 ("ID" / integer / FALSE) - this shifts rows
# 21..42 up to become rows 20..41.
$wsModel.Rows.Item(20).Delete()

# Insert a fresh row in front of what is now row 42 (formerly row 43),
# restoring the original row numbering for everything below it, and
# fill it in with the new "REGID" model entry.
$wsModel.Rows.Item(42).Insert()
$wsModel.Range("A42").Value = "REGID"
$wsModel.Range("B42").Value = "integer"
$wsModel.Range("C42").Value = $false

# Record the re-sorted range (A2:C50) as the sheet's remembered sort
# state, matching the tightened data block after the row swap above.
$sortObj = $wsModel.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($wsModel.Range("A2:A50"))
$sortObj.SetRange($wsModel.Range("A2:C50"))
$sortObj.Header = 0
$sortObj.Apply()

# --- view state (best effort) -----------------------------------------
$wsSurvey.Activate()
$wsSurvey.Range("F63").Select()

$wsModel.Activate()
$wsModel.Range("A1:C1048576").Select()

Write-Output "Renamed ID -> REGID (survey!F61, model rows) and refreshed model sortState."
